$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Detailed Procedure")

# Rows where column H currently holds "V3" and should become the text "9".
# A leading apostrophe forces Excel to keep the numeric-looking entry as
# text (quote-prefixed), matching how a user would type it in the UI,
# instead of letting it coerce to a number.
$rowsToNine = @(165, 166, 167, 177, 178, 179, 189, 190, 191)
foreach ($r in $rowsToNine) {
    $ws.Cells.Item($r, 8).Formula = "'9"
}

# Rows where column H currently holds "9" and should become the text "V4".
$rowsToV4 = @(168, 169, 170, 180, 181, 182, 192, 193, 194)
foreach ($r in $rowsToV4) {
    $ws.Cells.Item($r, 8).Formula = "V4"
}
